$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 108 (shifts existing rows 108-149 down to 109-150)
$ws.Rows("108:108").Insert()

# Populate the newly inserted row with the new price record
$ws.Range("A108").Value = 10
$ws.Range("B108").Value = "Vega Modelo de Temuco"
$ws.Range("C108").Value = "La Araucanía"
$ws.Range("D108").Value = 44510
$ws.Range("E108").Value = 9
$ws.Range("F108").Value = 100112013
$ws.Range("G108").Value = "Alcachofa"
$ws.Range("H108").Value = "Española"
$ws.Range("I108").Value = "Primera"
$ws.Range("J108").Value = 1400
$ws.Range("K108").Value = 400
$ws.Range("L108").Value = 500
$ws.Range("M108").Value = 457
$ws.Range("N108").Value = "$/unidad"
$ws.Range("O108").Value = "Región Metropolitana"
$ws.Range("P108").Value = 457
$ws.Range("Q108").Value = 1
$ws.Range("R108").Value = "Hortaliza"
